$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "329.41"
Set-TextValue "E2" "5.41%"
Set-TextValue "G2" "20"
Set-TextValue "D3" "40.41"
Set-TextValue "E3" "9.50%"
Set-TextValue "G3" "20"
Set-TextValue "D4" "5.794"
Set-TextValue "E4" "12.82%"
Set-TextValue "G4" "20"
Set-TextValue "D5" "0.08127"
Set-TextValue "E5" "3.51%"
Set-TextValue "G5" "20"
Set-TextValue "D6" "4.593"
Set-TextValue "E6" "4.47%"
Set-TextValue "G6" "20"
Set-TextValue "D7" "8.763"
Set-TextValue "E7" "4.41%"
Set-TextValue "G7" "20"
Set-TextValue "D8" "1.972"
Set-TextValue "E8" "4.78%"
Set-TextValue "G8" "20"
Set-TextValue "G9" "20"
Set-TextValue "D10" "0.9453"
Set-TextValue "E10" "1.61%"
Set-TextValue "G10" "20"
Set-TextValue "D11" "0.1319"
Set-TextValue "E11" "9.38%"
Set-TextValue "G11" "20"
Set-TextValue "D12" "0.1993"
Set-TextValue "E12" "4.24%"
Set-TextValue "G12" "20"
Set-TextValue "D13" "9.034"
Set-TextValue "E13" "37.69%"
Set-TextValue "G13" "20"
Set-TextValue "D14" "0.09307"
Set-TextValue "E14" "3.79%"
Set-TextValue "G14" "20"
Set-TextValue "D15" "0.03442"
Set-TextValue "E15" "4.61%"
Set-TextValue "G15" "20"
Set-TextValue "D16" "0.09608"
Set-TextValue "E16" "0.59%"
Set-TextValue "G16" "20"
Set-TextValue "D17" "0.001311"
Set-TextValue "E17" "-5.07%"
Set-TextValue "G17" "20"
Set-TextValue "D18" "0.006296"
Set-TextValue "E18" "6.55%"
Set-TextValue "G18" "20"
Set-TextValue "D19" "3.359"
Set-TextValue "E19" "-0.24%"
Set-TextValue "G19" "20"
Set-TextValue "D20" "0.3539"
Set-TextValue "E20" "2.31%"
Set-TextValue "G20" "20"
Set-TextValue "D21" "0.1405"
Set-TextValue "E21" "8.44%"
Set-TextValue "G21" "20"
Set-TextValue "D22" "0.2414"
Set-TextValue "E22" "5.01%"
Set-TextValue "G22" "20"
Set-TextValue "D23" "0.04429"
Set-TextValue "E23" "2.02%"
Set-TextValue "G23" "20"
Set-TextValue "D24" "0.001263"
Set-TextValue "E24" "5.66%"
Set-TextValue "G24" "20"
Set-TextValue "D25" "0.004385"
Set-TextValue "E25" "0.88%"
Set-TextValue "G25" "20"
Set-TextValue "E26" "-17.43%"
Set-TextValue "G26" "20"
Set-TextValue "D27" "0.0003999"
Set-TextValue "E27" "0.99%"
Set-TextValue "G27" "20"
Set-TextValue "G28" "20"
Set-TextValue "G29" "20"
Set-TextValue "G30" "20"
Set-TextValue "G31" "20"
Set-TextValue "G32" "20"
Set-TextValue "G33" "20"
Set-TextValue "G34" "20"
Set-TextValue "G35" "20"
Set-TextValue "G36" "20"
Set-TextValue "G37" "20"
Set-TextValue "G38" "20"
Set-TextValue "D39" "0.02466"
Set-TextValue "E39" "9.06%"
Set-TextValue "G39" "20"
Set-TextValue "D40" "0.05287"
Set-TextValue "E40" "3.46%"
Set-TextValue "G40" "20"
Set-TextValue "D41" "0.007484"
Set-TextValue "E41" "0.22%"
Set-TextValue "G41" "20"
Set-TextValue "D42" "0.1435"
Set-TextValue "E42" "3.78%"
Set-TextValue "G42" "20"
Set-TextValue "D43" "0.008934"
Set-TextValue "E43" "7.03%"
Set-TextValue "G43" "20"
Set-TextValue "D44" "0.002054"
Set-TextValue "E44" "3.35%"
Set-TextValue "G44" "20"
Set-TextValue "E45" "33.86%"
Set-TextValue "G45" "20"
Set-TextValue "D46" "0.00006893"
Set-TextValue "E46" "8.89%"
Set-TextValue "G46" "20"
Set-TextValue "E47" "0.75%"
Set-TextValue "G47" "20"
Set-TextValue "D48" "0.003503"
Set-TextValue "E48" "22.83%"
Set-TextValue "G48" "20"
Set-TextValue "D49" "0.001704"
Set-TextValue "E49" "1.37%"
Set-TextValue "G49" "20"
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.75%"
Set-TextValue "G50" "20"
Set-TextValue "D51" "0.0002005"
Set-TextValue "E51" "0.75%"
Set-TextValue "G51" "20"
